$wb = $excel.ActiveWorkbook

# --- Sheet "Hoja1": update the conversion note text in A1 ---
$wsHoja1 = $wb.Worksheets.Item("Hoja1")
$newLine = "`n"
$text = "Conversión del día 💰" + $newLine +
        "✅ Dólar paralelo: 68" + $newLine +
        $newLine +
        "Binance" + $newLine +
        "✅ 1000 Bs = 1.82 = 6763.65 pesos" + $newLine +
        "✅ 6763.65 pesos = 1.81 = 914.4 Bs" + $newLine +
        $newLine +
        "Promedio competencia" + $newLine +
        "✅ Tasa pesos: 20" + $newLine +
        "✅ Tasa Bs: 20" + $newLine +
        "✅ % Ganancia: 20%"
$wsHoja1.Range("A1").Value = $text

# --- Sheet "tasas": update N10, O10, N12, O12 ---
$wsTasas = $wb.Worksheets.Item("tasas")
$wsTasas.Range("N10").Value = 550
$wsTasas.Range("O10").Value = 3720.01
$wsTasas.Range("N12").Value = 3738
$wsTasas.Range("O12").Value = 505.35
